$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry match-specific data and get rotated between rows.
# (Column A = sequential index, C = Div, D = Date stay fixed per physical row.)
$cols = @(2) + (5..30)

function Swap-Rows([int]$r1, [int]$r2) {
    foreach ($col in $cols) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Simple swaps (pairs of rows whose data got exchanged)
Swap-Rows 59 60
Swap-Rows 113 114
Swap-Rows 184 185
Swap-Rows 191 192

# Row 178/179/180 form a 3-cycle: new178 = old180, new179 = old178, new180 = old179
$row178 = @{}
$row179 = @{}
$row180 = @{}
foreach ($col in $cols) {
    $row178[$col] = $ws.Cells.Item(178, $col).Value2
    $row179[$col] = $ws.Cells.Item(179, $col).Value2
    $row180[$col] = $ws.Cells.Item(180, $col).Value2
}
foreach ($col in $cols) {
    $ws.Cells.Item(178, $col).Value = $row180[$col]
    $ws.Cells.Item(179, $col).Value = $row178[$col]
    $ws.Cells.Item(180, $col).Value = $row179[$col]
}
